$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    "Toast Defense",
    "Tiny Toy Tanks",
    "BatuGame",
    "Phobos",
    "Aurora: A Child's Journey",
    "Verse Surf",
    "Super Buckyball Tournament Preseason",
    "Daikaiju Daikessen: Versus",
    "Our Life: Beginnings & Always",
    "Astria",
    "Disc Space",
    "Beat Me! - Puppetonia Tournament",
    "Tribal Wars",
    "GranAge",
    "Arc Apellago",
    "Arc Apellago",
    "Swing Dunk (Open Beta)",
    "Perfect Vermin",
    "Lawyer Guy: Defender of Justice",
    "How Stories Die",
    "Coloring Book for Adults",
    "Apex Legends™",
    "Jendo: Origins",
    "Space Commander: War and Trade",
    "Chess'Extra",
    "Sword of the Necromancer - Prologue",
    "SpellForce 3: Versus Edition",
    "Remnants",
    "Battle Star",
    "Cubiscape 2"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
